#
# Refresh the three KPI placeholder textboxes on slide 1:
#   - KPI1 ("TextBox 12", id 13) is repositioned/resized in place.
#   - KPI2 ("TextBox 21", id 22) and KPI3 ("TextBox 24", id 25) are replaced
#     by fresh copies of the KPI1 textbox (same run/paragraph formatting),
#     repositioned, retexted and renamed; the old shapes are removed.
#
# Point<->EMU note: Shape.Left/Top/Width/Height round-trip through a 32-bit
# float in this COM host (like real PowerPoint), so the literals below were
# chosen so that pt*12700 lands exactly on the target EMU value instead of
# the naive target_emu/12700 (which can be off by one EMU after the
# float32 round-trip).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- KPI1 ("TextBox 12") : reposition/resize only ----------------------
$kpi1 = $s.Shapes.Item("TextBox 12")
$kpi1.Left   = 11.531889763779528   # 146455 EMU
$kpi1.Top    = 47.90149606299212    # 608349 EMU
$kpi1.Width  = 171.72600393700787   # 2180920 EMU
$kpi1.Height = 29.081259842519685   # 369332 EMU

# ---- KPI2 ("TextBox 21" -> "TextBox 31") --------------------------------
$oldKpi2 = $s.Shapes.Item("TextBox 21")
$newKpi2 = $kpi1.Duplicate()
$newKpi2.Name = "TextBox 31"
$newKpi2.TextFrame.TextRange.Text = "KPI2"
$newKpi2.Left   = 12.248523622047243   # 155556 EMU
$newKpi2.Top    = 208.03049212598427   # 2641987 EMU
$newKpi2.Width  = 171.72600393700787   # 2180920 EMU
$newKpi2.Height = 29.081259842519685   # 369332 EMU
$oldKpi2.Delete()

# ---- KPI3 ("TextBox 24" -> "TextBox 32") --------------------------------
$oldKpi3 = $s.Shapes.Item("TextBox 24")
$newKpi3 = $kpi1.Duplicate()
$newKpi3.Name = "TextBox 32"
$newKpi3.TextFrame.TextRange.Text = "KPI3"
$newKpi3.Left   = 12.248523622047243   # 155556 EMU
$newKpi3.Top    = 366.15740157480315   # 4650199 EMU
$newKpi3.Width  = 171.72600393700787   # 2180920 EMU
$newKpi3.Height = 29.081259842519685   # 369332 EMU
$oldKpi3.Delete()
